# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 263
$ws1.Range("F5").Value = 3094
$ws1.Range("F6").Value = 2065
$ws1.Range("F9").Value = 1159
$ws1.Range("F10").Value = 211
$ws1.Range("F11").Value = 929
$ws1.Range("F12").Value = 80

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 263
$ws4.Range("F5").Value = 3094
$ws4.Range("F6").Value = 2065
$ws4.Range("F10").Value = 1159
$ws4.Range("F11").Value = 211
$ws4.Range("F12").Value = 929
$ws4.Range("F13").Value = 80
